$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @(2, 463.35, '(50.436, 30.528)'),
    @(3, 463.34, '(50.437, 30.527)'),
    @(4, 463.32, '(50.436, 30.527)'),
    @(5, 463.31, '(50.437, 30.528)'),
    @(6, 463.28, '(50.436, 30.529)'),
    @(7, 463.27, '(50.437, 30.526)'),
    @(8, 463.21, '(50.436, 30.526)'),
    @(9, 463.2, '(50.437, 30.529)'),
    @(10, 463.19, '(50.435, 30.528)'),
    @(11, 463.18, '(50.435, 30.529)'),
    @(12, 463.15, '(50.438, 30.527)'),
    @(13, 463.13, '(50.436, 30.53)'),
    @(14, 463.13, '(50.438, 30.526)'),
    @(15, 463.12, '(50.435, 30.527)'),
    @(16, 463.11, '(50.437, 30.525)'),
    @(17, 463.07, '(50.435, 30.53)'),
    @(18, 463.07, '(50.438, 30.528)'),
    @(19, 463.01, '(50.438, 30.525)'),
    @(20, 463, '(50.436, 30.525)'),
    @(21, 462.99, '(50.437, 30.53)'),
    @(22, 462.96, '(50.435, 30.526)'),
    @(23, 462.9, '(50.438, 30.529)'),
    @(24, 462.89, '(50.434, 30.529)'),
    @(25, 462.88, '(50.436, 30.531)'),
    @(26, 462.88, '(50.435, 30.531)'),
    @(27, 462.86, '(50.434, 30.528)'),
    @(28, 462.85, '(50.437, 30.524)'),
    @(29, 462.84, '(50.434, 30.53)'),
    @(30, 462.79, '(50.438, 30.524)'),
    @(31, 462.77, '(50.439, 30.526)'),
    @(32, 462.75, '(50.434, 30.527)'),
    @(33, 462.74, '(50.439, 30.527)'),
    @(34, 462.71, '(50.435, 30.525)'),
    @(35, 462.71, '(50.436, 30.524)'),
    @(36, 462.7, '(50.439, 30.525)'),
    @(37, 462.69, '(50.434, 30.531)'),
    @(38, 462.68, '(50.437, 30.531)'),
    @(39, 462.63, '(50.438, 30.53)'),
    @(40, 462.61, '(50.439, 30.528)'),
    @(41, 462.59, '(50.435, 30.532)'),
    @(42, 462.55, '(50.434, 30.526)'),
    @(43, 462.54, '(50.436, 30.532)'),
    @(44, 462.52, '(50.439, 30.524)'),
    @(45, 462.5, '(50.437, 30.523)'),
    @(46, 462.48, '(50.438, 30.523)'),
    @(47, 462.45, '(50.434, 30.532)'),
    @(48, 462.45, '(50.433, 30.529)'),
    @(49, 462.44, '(50.433, 30.53)'),
    @(50, 462.38, '(50.439, 30.529)'),
    @(51, 462.38, '(50.435, 30.524)'),
    @(52, 462.38, '(50.433, 30.528)'),
    @(53, 462.33, '(50.433, 30.531)'),
    @(54, 462.33, '(50.436, 30.523)'),
    @(55, 462.28, '(50.437, 30.532)'),
    @(56, 462.27, '(50.438, 30.531)'),
    @(57, 462.26, '(50.434, 30.525)'),
    @(58, 462.24, '(50.439, 30.523)'),
    @(59, 462.22, '(50.433, 30.527)'),
    @(60, 462.21, '(50.435, 30.533)'),
    @(61, 462.19, '(50.44, 30.526)'),
    @(62, 462.16, '(50.44, 30.525)'),
    @(63, 462.14, '(50.433, 30.532)'),
    @(64, 462.12, '(50.434, 30.533)'),
    @(65, 462.11, '(50.44, 30.527)'),
    @(66, 462.1, '(50.436, 30.533)'),
    @(67, 462.07, '(50.438, 30.522)'),
    @(68, 462.06, '(50.437, 30.522)'),
    @(69, 462.06, '(50.439, 30.53)'),
    @(70, 462.03, '(50.44, 30.524)'),
    @(71, 461.99, '(50.433, 30.526)'),
    @(72, 461.97, '(50.435, 30.523)'),
    @(73, 461.92, '(50.44, 30.528)'),
    @(74, 461.9, '(50.434, 30.524)'),
    @(75, 461.89, '(50.432, 30.53)'),
    @(76, 461.87, '(50.432, 30.529)'),
    @(77, 461.86, '(50.439, 30.522)'),
    @(78, 461.85, '(50.433, 30.533)'),
    @(79, 461.85, '(50.436, 30.522)'),
    @(80, 461.83, '(50.432, 30.531)'),
    @(81, 461.8, '(50.438, 30.532)'),
    @(82, 461.79, '(50.44, 30.523)'),
    @(83, 461.78, '(50.437, 30.533)'),
    @(84, 461.76, '(50.432, 30.528)'),
    @(85, 461.72, '(50.435, 30.534)'),
    @(86, 461.69, '(50.434, 30.534)'),
    @(87, 461.67, '(50.432, 30.532)'),
    @(88, 461.67, '(50.433, 30.525)'),
    @(89, 461.64, '(50.44, 30.529)'),
    @(90, 461.63, '(50.439, 30.531)'),
    @(91, 461.57, '(50.432, 30.527)'),
    @(92, 461.56, '(50.436, 30.534)'),
    @(93, 461.55, '(50.438, 30.521)'),
    @(94, 461.52, '(50.437, 30.521)'),
    @(95, 461.47, '(50.433, 30.534)'),
    @(96, 461.46, '(50.435, 30.522)'),
    @(97, 461.45, '(50.434, 30.523)'),
    @(98, 461.44, '(50.44, 30.522)'),
    @(99, 461.43, '(50.432, 30.533)'),
    @(100, 461.39, '(50.441, 30.525)'),
    @(101, 461.38, '(50.439, 30.521)'),
)

foreach ($item in $data) {
    $r = $item[0]
    $aVal = $item[1]
    $bVal = $item[2]
    $ws.Cells.Item($r, 1).Value = $aVal
    $ws.Cells.Item($r, 2).Value = $bVal
}

$wb.Save()